# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Row 16 and row 17 hold the two "Periodo Mora" / "Valor Mora" entries for
# the same worker. The database refresh re-orders them (2102 before 2103)
# -- swap the "Periodo Mora" (col E) and "Valor Mora" (col F) values between
# the two rows; the rest of the row (Tipo Doc, N Doc, Nombre, Salario) stays
# untouched since it's identical for both rows anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodoE16 = $ws.Range("E16").Value2
$periodoE17 = $ws.Range("E17").Value2
$valorF16 = $ws.Range("F16").Value2
$valorF17 = $ws.Range("F17").Value2

$ws.Range("E16").Value = $periodoE17
$ws.Range("E17").Value = $periodoE16

$ws.Range("F16").Value = $valorF17
$ws.Range("F17").Value = $valorF16
